# NB3_ear.xlsx -> rev3 fab
# Applies:
#  - BOM row for MK1 microphone: supplier changed from Mouser/410-ICS-43434 to LCSC/C5656610,
#    and unit price updated from 2.2 to 1.34
#  - Row heights normalized to 12.8 (was 14.35) for all used rows
#  - Selection / active cell moved to E13 (was A3)
#  - Page setup paper size changed to A4 (9) (was Letter/1)
#  - Header/Footer now use an explicit Arial,Regular font

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM update: MK1 microphone supplier + unit price (row 4) ---
$ws.Range("F4").Value = "LCSC"
$ws.Range("G4").Value = "C5656610"
$ws.Range("H4").Value = 1.34

# --- Row heights: rows 1-14 from 14.35 to 12.8 ---
for ($r = 1; $r -le 14; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.8
}

# --- Selection / view state ---
$ws.Range("E13").Select()

# --- Page setup ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.CenterHeader = '&"Arial,Regular"&A'
$ps.CenterFooter = '&"Arial,Regular"Page &P'
